$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "29.524.16"
Set-TextCell "E2" "  -0.66%  "

Set-TextCell "D3" "1.850.88"
Set-TextCell "E3" "  -0.16%  "

Set-TextCell "E4" "  +0.01%  "

Set-TextCell "D5" "243.06"
Set-TextCell "E5" "  -0.42%  "

Set-TextCell "D6" "0.6364"
Set-TextCell "E6" "  -0.57%  "

Set-TextCell "E7" "  +0.05%  "

Set-TextCell "D8" "48.01"
Set-TextCell "E8" "  +1.58%  "

Set-TextCell "D9" "0.07574"
Set-TextCell "E9" "  +1.20%  "

Set-TextCell "D10" "0.2999"
Set-TextCell "E10" "  +0.48%  "

Set-TextCell "D11" "24.23"
Set-TextCell "E11" "  -0.63%  "

Set-TextCell "D12" "0.07695"
Set-TextCell "E12" "  +0.62%  "

Set-TextCell "D13" "1.874.90"
Set-TextCell "E13" "  +0.87%  "

Set-TextCell "D14" "5.037"
Set-TextCell "E14" "  -0.13%  "

Set-TextCell "D15" "0.6878"
Set-TextCell "E15" "  -0.03%  "

Set-TextCell "D16" "84.01"
Set-TextCell "E16" "  +0.23%  "

Set-TextCell "D17" "0.000009822"
Set-TextCell "E17" "  +3.49%  "

Set-TextCell "D18" "2.110.88"
Set-TextCell "E18" "  +0.34%  "

Set-TextCell "D19" "6.276"
Set-TextCell "E19" "  +3.67%  "

Set-TextCell "D20" "29.562.18"
Set-TextCell "E20" "  -0.55%  "

Set-TextCell "D21" "237.23"
Set-TextCell "E21" "  +0.66%  "

Set-TextCell "E22" "  -0.67%  "

Set-TextCell "E23" "  +0.06%  "

Set-TextCell "D24" "7.621"
Set-TextCell "E24" "  +2.81%  "

Set-TextCell "E25" "  -0.01%  "

Set-TextCell "D26" "156.50"
Set-TextCell "E26" "  -1.18%  "

Set-TextCell "D27" "0.1395"
Set-TextCell "E27" "  -1.61%  "

Set-TextCell "D28" "8.455"
Set-TextCell "E28" "  -0.37%  "

Set-TextCell "E29" "  -0.72%  "

Set-TextCell "E30" "  -0.30%  "

Set-TextCell "D31" "0.05891"
Set-TextCell "E31" "  -6.19%  "

Set-TextCell "D32" "1.277"
Set-TextCell "E32" "  -0.01%  "

Set-TextCell "D33" "4.130"
Set-TextCell "E33" "  -0.40%  "

Set-TextCell "D34" "4.070"
Set-TextCell "E34" "  -0.50%  "

Set-TextCell "D35" "1.900"
Set-TextCell "E35" "  -0.09%  "

Set-TextCell "E36" "  +0.34%  "

Set-TextCell "D37" "0.7200"
Set-TextCell "E37" "  -1.21%  "

Set-TextCell "D38" "2.599"
Set-TextCell "E38" "  -0.24%  "

Set-TextCell "D39" "2.809"
Set-TextCell "E39" "  -1.22%  "

Set-TextCell "D40" "1.228.96"
Set-TextCell "E40" "  +2.20%  "

Set-TextCell "D41" "0.01779"
Set-TextCell "E41" "  -0.26%  "

Set-TextCell "D42" "0.9130"
Set-TextCell "E42" "  -0.99%  "

Set-TextCell "D43" "6.124"
Set-TextCell "E43" "  -0.49%  "

Set-TextCell "D44" "0.9997"

Set-TextCell "D45" "2.020.52"
Set-TextCell "E45" "  +0.31%  "

Set-TextCell "D46" "101.94"
Set-TextCell "E46" "  -0.02%  "

Set-TextCell "D47" "67.48"
Set-TextCell "E47" "  +2.17%  "

Set-TextCell "D48" "7.397"
Set-TextCell "E48" "  +10.67%  "

Set-TextCell "D49" "0.4046"
Set-TextCell "E49" "  -0.39%  "

# Rows 50 and 51 swap (EnergySwap <-> BabyDogeCoin) plus updated price/volume
Set-TextCell "B50" "BabyDogeCoin"
Set-TextCell "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D50" "0.00000000117"
Set-TextCell "E50" "  -1.37%  "
Set-TextCell "B51" "EnergySwap"
Set-TextCell "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D51" "9.144"
Set-TextCell "E51" "  -0.76%  "
